$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their literal string values (no numeric/percent auto-conversion)
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Column D (Price) updates
$ws.Range("D2").Value = "260.91"
$ws.Range("D3").Value = "26.84"
$ws.Range("D4").Value = "4.705"
$ws.Range("D5").Value = "0.06221"
$ws.Range("D6").Value = "6.741"
$ws.Range("D7").Value = "0.8496"
$ws.Range("D8").Value = "0.9136"
$ws.Range("D9").Value = "0.1402"
$ws.Range("D10").Value = "0.04948"
$ws.Range("D11").Value = "0.07086"
$ws.Range("D12").Value = "0.03084"
$ws.Range("D13").Value = "0.09050"
$ws.Range("D14").Value = "0.001533"
$ws.Range("D15").Value = "0.0006165"
$ws.Range("D16").Value = "0.005949"
$ws.Range("D17").Value = "3.447"
$ws.Range("D19").Value = "2.168"
$ws.Range("D21").Value = "0.1311"
$ws.Range("D22").Value = "4.092"
$ws.Range("D23").Value = "0.04245"
$ws.Range("D24").Value = "0.001204"
$ws.Range("D40").Value = "0.03961"
$ws.Range("D41").Value = "0.1113"
$ws.Range("D42").Value = "0.004135"
$ws.Range("D43").Value = "0.002101"
$ws.Range("D44").Value = "0.01333"
$ws.Range("D45").Value = "0.00005162"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D48").Value = "0.2395"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.0002001"

# Column E (Volume(1h)) updates
$ws.Range("E3").Value = "-1.44%"
$ws.Range("E4").Value = "-0.44%"
$ws.Range("E5").Value = "2.29%"
$ws.Range("E6").Value = "1.08%"
$ws.Range("E7").Value = "0.50%"
$ws.Range("E8").Value = "-1.00%"
$ws.Range("E9").Value = "-0.43%"
$ws.Range("E10").Value = "-0.67%"
$ws.Range("E11").Value = "-0.17%"
$ws.Range("E12").Value = "-1.33%"
$ws.Range("E13").Value = "-0.26%"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("E15").Value = "1.16%"
$ws.Range("E16").Value = "-3.82%"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("E18").Value = "0.85%"
$ws.Range("E19").Value = "0.01%"
$ws.Range("E21").Value = "1.08%"
$ws.Range("E22").Value = "-0.23%"
$ws.Range("E23").Value = "0.09%"
$ws.Range("E24").Value = "-1.45%"
$ws.Range("E25").Value = "4.26%"
$ws.Range("E26").Value = "0.03%"
$ws.Range("E27").Value = "4.11%"
$ws.Range("E40").Value = "2.18%"
$ws.Range("E41").Value = "-0.01%"
$ws.Range("E42").Value = "0.67%"
$ws.Range("E43").Value = "-4.87%"
$ws.Range("E44").Value = "-18.45%"
$ws.Range("E45").Value = "-2.96%"
$ws.Range("E46").Value = "0.02%"
$ws.Range("E48").Value = "76.34%"
$ws.Range("E49").Value = "0.02%"
$ws.Range("E50").Value = "0.02%"

# Column G (Hora) updates: every data row 2-51 changes from 7 to 8
$ws.Range("G2:G51").Value = "8"

